# menu item deletion when branch is closed
# Delete the data row for the "NTU" branch (row 2) — Excel shifts all
# subsequent rows up by one automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
